$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The English text in row 74 ("Initial & Final Surveillance Diagnosis") was not
# translated correctly because of the literal "&" character. Fix this by:
#  1) Inserting a new row 75 that holds the corrected English text
#     ("Initial and Final Surveillance Diagnosis") together with the Lao
#     translation that used to sit on row 74 (column B).
#  2) Marking the original row 74 (with the literal "&") as "TBT" (to be
#     translated) in column B, since it is effectively a duplicate/broken key.
#  3) Fixing the other occurrence of the same "&" bug further down
#     ("Susceptible & Intermediate..."), which after the insertion above now
#     lives on row 145.

$oldB74 = $ws.Range("B74").Value()

$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value = $oldB74

$ws.Range("B74").Value = "TBT"

$ws.Range("A145").Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
